$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Customer_name column (column B) entirely - header and all values
$ws.Range("B1:B11").ClearContents()

# Update selection to match the new active cell
$ws.Range("B8").Select()
